$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: direct assignment ---
$ws.Range("D2").Value = "35.432.24"
$ws.Range("E2").Value = "  +2.80%  "
$ws.Range("D3").Value = "1.850.36"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +9.71%  "
$ws.Range("E9").Value = "  +6.39%  "
$ws.Range("E10").Value = "  +3.15%  "
$ws.Range("E11").Value = "  +3.91%  "
$ws.Range("D12").Value = "2.119.29"
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").Value = "1.850.46"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("E15").Value = "  +6.81%  "
$ws.Range("D17").Value = "35.445.19"
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("D20").Value = "0.0₃0798"
$ws.Range("E20").Value = "  +3.88%  "
$ws.Range("E21").Value = "  +8.67%  "
$ws.Range("E22").Value = "  +13.22%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("E29").Value = "  +12.88%  "
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "3.289.66"
$ws.Range("E31").Value = "  +35.39%  "
$ws.Range("E32").Value = "  +5.95%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E34").Value = "  +6.08%  "
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("E36").Value = "  +20.27%  "
$ws.Range("E37").Value = "  +6.71%  "
$ws.Range("D38").Value = "1.354.69"
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("E39").Value = "  +2.96%  "
$ws.Range("E40").Value = "  +5.56%  "
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("E42").Value = "  +5.62%  "
$ws.Range("E43").Value = "  +4.23%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E45").Value = "  +7.50%  "
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("E48").Value = "  +7.47%  "
$ws.Range("D49").Value = "2.017.67"
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("E51").Value = "  +1.64%  "

# --- Numeric-looking values that must remain TEXT: force text, then strip the
#     resulting "Text" number-format style by re-pasting the (unstyled) format
#     from an untouched reference cell, so no new style is left attached. ---
$fmtSource = $ws.Range("A1")
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.62"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.609"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.91"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.42"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.87"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.38"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.06"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.64"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.86"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.91"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.72"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.92"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.05"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "98.12"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.682"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.26"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.48"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.62"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0520"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.19"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.57"

# Restore original (default) formatting on those cells so no stray style remains
# (multi-area ranges only reliably paste into the first area, so loop per-cell)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null
$ws.Range("D21").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("D33").PasteSpecial(-4122) | Out-Null
$ws.Range("D34").PasteSpecial(-4122) | Out-Null
$ws.Range("D36").PasteSpecial(-4122) | Out-Null
$ws.Range("D37").PasteSpecial(-4122) | Out-Null
$ws.Range("D42").PasteSpecial(-4122) | Out-Null
$ws.Range("D43").PasteSpecial(-4122) | Out-Null
$ws.Range("D44").PasteSpecial(-4122) | Out-Null
$ws.Range("D45").PasteSpecial(-4122) | Out-Null
$ws.Range("D47").PasteSpecial(-4122) | Out-Null
$ws.Range("D48").PasteSpecial(-4122) | Out-Null
$ws.Range("D51").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false